# Apply FlashScore odds updates for 2025-04-17 (commit: "Atualizando o arquivo XLSX")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 2.2
$ws.Range("I2").Value = 3.3
$ws.Range("T2").Value = 7
$ws.Range("U2").Value = 10
$ws.Range("W2").Value = 21
$ws.Range("Z2").Value = 8.5
$ws.Range("AE2").Value = 8.5
$ws.Range("AF2").Value = 15
$ws.Range("AH2").Value = 34
# Row 5
$ws.Range("J5").Value = 1.05
$ws.Range("L5").Value = 1.33
# Row 6
$ws.Range("J6").Value = 1.05
$ws.Range("L6").Value = 1.3
# Row 7
$ws.Range("J7").Value = 1.07
$ws.Range("L7").Value = 1.41
$ws.Range("M7").Value = 2.62
# Row 8
$ws.Range("G8").Value = 1.7
$ws.Range("I8").Value = 5.75
$ws.Range("J8").Value = 1.05
$ws.Range("K8").Value = 8
$ws.Range("L8").Value = 1.3
$ws.Range("AD8").Value = 1000
$ws.Range("AF8").Value = 26
$ws.Range("AH8").Value = 51
# Row 9
$ws.Range("J9").Value = 1.03
$ws.Range("K9").Value = 15
$ws.Range("L9").Value = 1.2
$ws.Range("M9").Value = 4.33
# Row 11
$ws.Range("N11").Value = 2.35
$ws.Range("O11").Value = 1.57
# Row 12
$ws.Range("G12").Value = 1.52
$ws.Range("H12").Value = 3.75
$ws.Range("L12").Value = 1.39
$ws.Range("M12").Value = 2.55
$ws.Range("N12").Value = 2.15
$ws.Range("O12").Value = 1.55
$ws.Range("P12").Value = 1.47
$ws.Range("Q12").Value = 2.32
$ws.Range("R12").Value = 2.27
$ws.Range("S12").Value = 1.5
$ws.Range("T12").Value = 5
$ws.Range("U12").Value = 5.8
$ws.Range("V12").Value = 9
$ws.Range("X12").Value = 15
$ws.Range("Y12").Value = 45
$ws.Range("Z12").Value = 7.7
$ws.Range("AA12").Value = 7.6
$ws.Range("AB12").Value = 26
$ws.Range("AJ12").Value = 110
# Row 14
$ws.Range("G14").Value = 2.18
$ws.Range("H14").Value = 3.3
$ws.Range("I14").Value = 3.05
$ws.Range("M14").Value = 2.95
$ws.Range("N14").Value = 1.87
$ws.Range("R14").Value = 1.72
$ws.Range("S14").Value = 1.9
$ws.Range("T14").Value = 7.5
$ws.Range("U14").Value = 10.25
$ws.Range("V14").Value = 9
$ws.Range("W14").Value = 21
$ws.Range("X14").Value = 18
$ws.Range("Z14").Value = 9.75
$ws.Range("AA14").Value = 6.4
$ws.Range("AE14").Value = 9.75
$ws.Range("AF14").Value = 16
$ws.Range("AG14").Value = 11
$ws.Range("AH14").Value = 40
$ws.Range("AI14").Value = 26
$ws.Range("AJ14").Value = 35
# Row 16
$ws.Range("N16").Value = 1.65
$ws.Range("O16").Value = 2.2
# Row 17
$ws.Range("G17").Value = 1.3
$ws.Range("H17").Value = 5.5
$ws.Range("K17").Value = 17
$ws.Range("P17").Value = 1.29
$ws.Range("Q17").Value = 3.5
$ws.Range("R17").Value = 1.95
$ws.Range("S17").Value = 1.8
$ws.Range("T17").Value = 8
$ws.Range("Y17").Value = 26
$ws.Range("AE17").Value = 23
$ws.Range("AG17").Value = 26
# Row 18
$ws.Range("J18").Value = 1.06
$ws.Range("K18").Value = 10
# Row 19
$ws.Range("L19").Value = 1.25
$ws.Range("M19").Value = 3.75
$ws.Range("N19").Value = 1.8
$ws.Range("O19").Value = 2
# Row 20
$ws.Range("K20").Value = 15
$ws.Range("N20").Value = 1.29
$ws.Range("O20").Value = 3.5
# Row 21
$ws.Range("G21").Value = 1.48
$ws.Range("H21").Value = 4.2
$ws.Range("I21").Value = 6
$ws.Range("L21").Value = 1.17
$ws.Range("M21").Value = 4.5
$ws.Range("P21").Value = 1.3
$ws.Range("Q21").Value = 3.4
$ws.Range("T21").Value = 8.5
$ws.Range("W21").Value = 11
$ws.Range("Z21").Value = 15
$ws.Range("AA21").Value = 8.5
$ws.Range("AE21").Value = 19
$ws.Range("AF21").Value = 34
$ws.Range("AG21").Value = 19
# Row 22
$ws.Range("G22").Value = 3.9
$ws.Range("I22").Value = 1.75
$ws.Range("J22").Value = 23
$ws.Range("K22").Value = 1.02
$ws.Range("L22").Value = 1.07
$ws.Range("P22").Value = 1.22
$ws.Range("Q22").Value = 4
$ws.Range("U22").Value = 26
$ws.Range("AB22").Value = 12
$ws.Range("AD22").Value = 81
$ws.Range("AE22").Value = 13
